# Fixed Bento 80 Test scripts
# Append "ORDER BY ... LIMIT 100" clauses to the Neo4j queries stored in
# column B (the "query" column) of the startup sheet, for the CasesTab,
# SamplesTab and FilesTab rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Row 3 - SamplesTab query: append a new ORDER BY / LIMIT line.
$samplesQuery = $ws.Range("B3").Value2
$ws.Range("B3").Value = $samplesQuery + $nl + "order By samp.sample_id ASC LIMIT 100"

# Row 2 - CasesTab query: append a new ORDER BY / LIMIT line.
$casesQuery = $ws.Range("B2").Value2
$ws.Range("B2").Value = $casesQuery + $nl + " order By ss.study_subject_id ASC LIMIT 100"

# Row 4 - FilesTab query: replace the old trailing "order by f.file_name"
# clause with the new capitalized ORDER BY / LIMIT clause.
$filesQuery = $ws.Range("B4").Value2
$oldTail = "    order by f.file_name"
if ($filesQuery.EndsWith($oldTail)) {
    $filesQuery = $filesQuery.Substring(0, $filesQuery.Length - $oldTail.Length)
}
$ws.Range("B4").Value = $filesQuery + " order By f.file_name ASC LIMIT 100"

# The extra wrapped line in rows 2 and 3 makes those rows taller; row 4 was
# already at Excel's row-height cap so it is unaffected.
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360

# Scroll the view up slightly, matching the saved view position.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
